# edit.ps1 -- transform the "Cybersecurity" essay into "The Art of Expression" essay
# per the target commit diff.

$d = $word.ActiveDocument

# Manual line-break character, as produced by Shift+Enter (renders as <w:br/>)
$brk = [char]11

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "WARNING: not found -> $old"
    }
}

# ---- Title / byline / contact block ---------------------------------------

Replace-Text "Cybersecurity: Protecting the Digital Landscape" "The Art of Expression: Exploring Human Creativity Through the Arts"

Replace-Text "Daniella Ferguson" "Ethan Jones"

Replace-Text "daniella" "ejones(at)hs"

Replace-Text "ferguson@xyzmail.org" "edu"

# ---- Body paragraph 1 (intro) ----------------------------------------------

Replace-Text "At the heart of our interconnected world lies cybersecurity, a domain of immense significance in the digital era" "In the annals of human civilization, there lies a realm where imagination, emotion, and expression converge to weave a tapestry of beauty and meaning: the realm of the arts"

Replace-Text " From the intricate cyber networks that underpin critical infrastructure to the personal data entrusted to digital platforms, safeguarding these realms is of paramount importance" " Across cultures and time, humans have harnessed the power of art to communicate, explore ideas, and transcend the boundaries of reality"

Replace-Text " Cybersecurity is not merely a reactive response to threats, but a proactive endeavor that compels us to stay vigilant against a constantly evolving landscape of risks" " Art, in its myriad forms, invites us to delve into the depths of the human experience, unlocking hidden truths, stirring emotions, and inspiring thought"

# merges "...essential services" + "." + " Nonetheless...placed upon it" into one sentence
Replace-Text " The consequences of neglecting cybersecurity can be dire, resulting in data breaches, financial losses, and disruptions in essential services. Nonetheless, it is through collaboration, innovation, and education that we can bolster the resilience of our digital infrastructure, ensuring its continued integrity and safeguarding the trust placed upon it" " In this essay, we shall embark on a journey to unveil the significance of the arts in shaping human lives and societies"

# ---- Body paragraph 1, second "section" (after the blank double-break) ----

Replace-Text "Cybersecurity demands a multi-pronged approach that encompasses technological advancements, legal frameworks, and societal awareness" "Firstly, art serves as a universal language that transcends cultural and linguistic barriers"

Replace-Text " The rapid pace of technological progress necessitates the continual development of security solutions, encryption techniques, and threat intelligence systems" " It possesses the remarkable ability to bridge gaps, connect people from diverse backgrounds, and foster a sense of shared humanity"

Replace-Text " Legal frameworks must keep pace with these developments, ensuring appropriate regulations and penalties for cybercrimes" " A melody can evoke joy or sorrow in the hearts of listeners across the globe, a painting can transport viewers to distant lands or forgotten eras, and a dance can communicate stories that words cannot"

# merges "...software updates" + "." + " Only through...interconnected world" into one sentence
Replace-Text " Furthermore, cultivating a security-conscious society is crucial, where individuals embrace best practices in digital hygiene, recognizing the importance of strong passwords, multi-factor authentication, and software updates. Only through such collective efforts can we mitigate the risks posed by cyber threats, securing the digital assets that underpin our interconnected world" " Through art, we find a common ground where differences dissolve, and unity prevails"

# ---- Body paragraph 1, third "section" ------------------------------------

Replace-Text "The interdependence of our digital infrastructure mandates international collaboration" "Furthermore, art is a potent force for self-expression"

Replace-Text " The cross-border nature of cyberspace means that threats transcend national boundaries, necessitating a united front against cybercrime" " It provides individuals with an outlet to articulate their innermost thoughts, feelings, and experiences"

Replace-Text " Information sharing and coordinated responses are essential to deterring and neutralizing sophisticated attacks that span multiple jurisdictions" " When words fail to capture the complexities of the human psyche, art steps in as a medium of expression that allows individuals to communicate their unique perspectives"

# drops the " Moreover, global " lead-in and rewrites the remainder of the sentence
Replace-Text " Moreover, global cooperation can foster the development of shared standards and protocols, harmonizing approaches to cybersecurity and facilitating collective action against emerging threats" " Whether it is through the brushstrokes of a painter, the melodies of a musician, or the movements of a dancer, art enables individuals to share their stories with the world, thereby validating their experiences and fostering a sense of empathy and connection among humanity"

# final sentence of the old third section is replaced by a double line-break plus
# an entire new fourth "section" of content
$newFourthSection = "$brk$brk" + "Finally, art possesses the transformative power to inspire positive change in individuals and societies. From the works of great thinkers and revolutionary leaders to the humble scribbles of a child, art has the capacity to ignite change, challenge societal norms, and promote social justice. Art can raise awareness, spark conversations, and mobilize individuals to action. It can question authority, expose oppression, and provide a platform for the marginalized to voice their concerns. Throughout history, art has been an instrumental force in shaping political, social, and cultural movements, contributing to the advancement of human rights, equality, and justice"

Replace-Text " By pooling our resources and expertise, we can create a more secure and stable cyberspace that benefits all nations and individuals alike" $newFourthSection

# ---- Summary paragraph ------------------------------------------------------

Replace-Text "Cybersecurity stands as the cornerstone of our digital age, safeguarding the integrity and confidentiality of our interconnected world" "In conclusion, the arts play a vital role in human lives and societies, transcending cultural and linguistic barriers, providing a medium for self-expression, and inspiring positive change"

Replace-Text " It encompasses a holistic approach that demands technological advancements, robust legal frameworks, and a security-conscious society" " Art connects people, fosters empathy, and allows individuals to explore the depths of their own emotions and experiences"

Replace-Text " International collaboration is paramount, fostering collective responses to global threats and promoting harmonized standards" " Through art, we learn about ourselves and the world around us, broadening our perspectives, challenging our assumptions, and cultivating a greater appreciation for beauty, creativity, and human ingenuity"

Replace-Text " By embracing innovation, education, and global cooperation, we can fortify our digital infrastructure, ensuring its resilience and protecting the trust placed upon it" " Art is not merely a luxury, but an essential component of what makes us human, reminding us of our capacity for imagination, empathy, and the pursuit of a more meaningful and fulfilling existence"

# ---- Structural: extra empty paragraph at the very end of the document body ----

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
